$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$ws = $wb.Worksheets.Item("Metadata")

# Version bump
$ws.Range("B3").Value = "6.0.0"

# Date update
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was empty)
$ws.Range("B9").Value = "Alvearie Team"

# Replace Contact/No display row with Jurisdiction/United States of America
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Remove the now-duplicate Contact row entirely (row 11), shifting rows 12-21 up
$ws.Rows.Item(11).Delete()

# --- Elements sheet updates ---
$ws2 = $wb.Worksheets.Item("Elements")

# Root Extension element (row 2): give it the specific Short/Definition text
$ws2.Range("K2").Value = "Nonstandard Language"
$ws2.Range("L2").Value = "Code for the language of the person"
